$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Collection_SC")

# Insert a new, blank column at K (shifts K:AJ -> L:AK for all rows)
$ws.Columns("K").Insert()

# Label the new column header
$ws.Range("K1").Value = "categories"

# The values that used to live under the "scenario" header (column J, now
# the original data) actually belong under the new "categories" column for
# the data rows. Move them over and clear the old cells.
$ws.Range("K2:K6").Value2 = $ws.Range("J2:J6").Value2
$ws.Range("J2:J6").ClearContents()

# Match the final cursor position left behind by the editing session
[void]$ws.Range("L14").Select()

Write-Output "done"
